$d = $word.ActiveDocument

$pairs = @(
    @("49÷8=6, 1", "80÷2=40, 0"),
    @("27÷7=3, 6", "14÷2=7, 0"),
    @("79÷6=13, 1", "61÷4=15, 1"),
    @("51÷3=17, 0", "34÷8=4, 2"),
    @("89÷5=17, 4", "28÷6=4, 4"),
    @("51÷4=12, 3", "59÷5=11, 4"),
    @("10÷7=1, 3", "89÷4=22, 1"),
    @("19÷3=6, 1", "40÷9=4, 4"),
    @("57÷3=19, 0", "92÷5=18, 2"),
    @("55÷2=27, 1", "81÷2=40, 1"),
    @("50÷7=7, 1", "70÷3=23, 1"),
    @("48÷2=24, 0", "82÷6=13, 4"),
    @("21÷2=10, 1", "95÷9=10, 5"),
    @("85÷2=42, 1", "80÷9=8, 8"),
    @("19÷8=2, 3", "17÷3=5, 2"),
    @("16÷6=2, 4", "33÷5=6, 3"),
    @("75÷2=37, 1", "59÷8=7, 3"),
    @("76÷2=38, 0", "76÷4=19, 0"),
    @("64÷4=16, 0", "60÷9=6, 6"),
    @("88÷9=9, 7", "24÷7=3, 3"),
    @("75÷7=10, 5", "13÷8=1, 5"),
    @("55÷6=9, 1", "26÷7=3, 5"),
    @("20÷2=10, 0", "42÷7=6, 0"),
    @("31÷6=5, 1", "87÷9=9, 6"),
    @("59÷7=8, 3", "84÷8=10, 4")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
